# Update cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the original inline-string cell type) instead of auto-converting
# numeric-looking strings (e.g. "318.13") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'28.343.30"
$ws.Range("E2").Value = "'  +5.39%  "
$ws.Range("D3").Value = "'1.810.85"
$ws.Range("E3").Value = "'  +4.55%  "
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("D5").Value = "'318.13"
$ws.Range("E5").Value = "'  +2.69%  "
$ws.Range("E6").Value = "'  +0.19%  "
$ws.Range("D7").Value = "'0.5716"
$ws.Range("E7").Value = "'  +14.58%  "
$ws.Range("D8").Value = "'0.3886"
$ws.Range("E8").Value = "'  +11.01%  "
$ws.Range("D9").Value = "'0.07595"
$ws.Range("E9").Value = "'  +4.78%  "
$ws.Range("D10").Value = "'42.94"
$ws.Range("E10").Value = "'  -0.04%  "
$ws.Range("D11").Value = "'1.140"
$ws.Range("E11").Value = "'  +8.05%  "
$ws.Range("E12").Value = "'  +0.17%  "
$ws.Range("D13").Value = "'21.19"
$ws.Range("E13").Value = "'  +6.32%  "
$ws.Range("D14").Value = "'6.267"
$ws.Range("E14").Value = "'  +6.59%  "
$ws.Range("D15").Value = "'1.810.23"
$ws.Range("E15").Value = "'  +4.77%  "
$ws.Range("D16").Value = "'7.277"
$ws.Range("E16").Value = "'  +6.79%  "
$ws.Range("D17").Value = "'92.02"
$ws.Range("E17").Value = "'  +5.97%  "
$ws.Range("E18").Value = "'  +3.70%  "
$ws.Range("D19").Value = "'0.06484"
$ws.Range("E19").Value = "'  +1.46%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "'  +0.17%  "
$ws.Range("D21").Value = "'17.31"
$ws.Range("E21").Value = "'  +4.53%  "
$ws.Range("D22").Value = "'6.006"
$ws.Range("D23").Value = "'28.357.60"
$ws.Range("E23").Value = "'  +5.20%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("E24").Value = "'  +3.14%  "
$ws.Range("D25").Value = "'2.141"
$ws.Range("E25").Value = "'  +4.02%  "
$ws.Range("E26").Value = "'  +18.21%  "
$ws.Range("D27").Value = "'158.20"
$ws.Range("E27").Value = "'  +2.84%  "
$ws.Range("D28").Value = "'20.79"
$ws.Range("E28").Value = "'  +4.35%  "
$ws.Range("D29").Value = "'2.020.09"
$ws.Range("E29").Value = "'  +4.85%  "
$ws.Range("D30").Value = "'124.22"
$ws.Range("E30").Value = "'  +3.43%  "
$ws.Range("D31").Value = "'1.163"
$ws.Range("E31").Value = "'  +11.14%  "
$ws.Range("D32").Value = "'0.1072"
$ws.Range("E32").Value = "'  +13.96%  "
$ws.Range("D33").Value = "'5.793"
$ws.Range("E33").Value = "'  +7.56%  "
$ws.Range("D34").Value = "'3.632"
$ws.Range("E34").Value = "'  +1.57%  "
$ws.Range("D35").Value = "'0.2227"
$ws.Range("E35").Value = "'  +12.13%  "
$ws.Range("D36").Value = "'8.985"
$ws.Range("E36").Value = "'  +20.94%  "
$ws.Range("D37").Value = "'0.02322"
$ws.Range("E37").Value = "'  +6.31%  "
$ws.Range("D38").Value = "'11.69"
$ws.Range("E38").Value = "'  +6.34%  "
$ws.Range("D39").Value = "'0.06128"
$ws.Range("E39").Value = "'  +3.64%  "
$ws.Range("D40").Value = "'0.6401"
$ws.Range("E40").Value = "'  +6.49%  "
$ws.Range("D41").Value = "'5.037"
$ws.Range("D42").Value = "'1.163"
$ws.Range("E42").Value = "'  +4.45%  "
$ws.Range("E43").Value = "'  +0.16%  "
$ws.Range("D44").Value = "'1.379"
$ws.Range("E44").Value = "'  -3.34%  "
$ws.Range("D45").Value = "'13.42"
$ws.Range("E45").Value = "'  +5.07%  "
$ws.Range("D46").Value = "'0.6012"
$ws.Range("E46").Value = "'  +7.08%  "
$ws.Range("D47").Value = "'3.694"
$ws.Range("E47").Value = "'  +3.39%  "
$ws.Range("D48").Value = "'122.75"
$ws.Range("E48").Value = "'  +2.62%  "
$ws.Range("D49").Value = "'1.953"
$ws.Range("E49").Value = "'  +6.17%  "
$ws.Range("D50").Value = "'1.148"
$ws.Range("E50").Value = "'  +4.58%  "
$ws.Range("E51").Value = "'  +3.27%  "
